$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for the season-record columns (Wins / Losses / Ties),
# appended right after the existing "Unnamed: 28" column (AC).
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the bold/bordered/centered header style used by the rest of row 1
# (copy format only, so the new cells reuse the same style as A1 instead of
# creating a near-duplicate style entry).
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Every player row (2-43) gets the same team season record: 88 wins,
# 74 losses, 1 tie.
for ($row = 2; $row -le 43; $row++) {
    $ws.Cells.Item($row, 30).Value = 88
    $ws.Cells.Item($row, 31).Value = 74
    $ws.Cells.Item($row, 32).Value = 1
}
